$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the data values (swap the -1 / 3 pair between B8 and B13)
$ws.Range("B8").Value = 3
$ws.Range("B13").Value = -1

# Update the current selection to match the new active cell
$ws.Range("B8").Select()
